$wb = $excel.ActiveWorkbook

$wsW22 = $wb.Worksheets.Item("W22")
$wsTopics = $wb.Worksheets.Item("Topics")

# Add new Management Process entry on the W22 schedule sheet
$wsW22.Range("D4").Value = "Management Process"
$wsW22.Range("E4").Value = "ManagementProcess"

# Move the active selection on W22 to D13
$wsW22.Range("D13").Select()

# Remove the now-duplicate first row's content from the Topics sheet
# (row numbering below it is unaffected - only the A1 cell content goes away)
$wsTopics.Range("A1").ClearContents()
